# Sign-up data providers updated:
#  - "VisaAdyen" sheet renamed to "Username" and rebuilt with username test-data rows
#  - "3DSUser" sheet renamed to "Password" and rebuilt with matching rows
#  - "FullName" sheet loses the active-tab state; "Username" becomes the active tab

$wb = $excel.ActiveWorkbook

$wsUsername = $wb.Worksheets.Item("VisaAdyen")
$wsPassword = $wb.Worksheets.Item("3DSUser")
$wsFullName = $wb.Worksheets.Item("FullName")

# --- Rename sheets ------------------------------------------------------
$wsUsername.Name = "Username"
$wsPassword.Name = "Password"

# --- Wipe old contents/formatting of both rebuilt sheets ----------------
$wsUsername.Cells.Clear()
$wsPassword.Cells.Clear()

# --- Column widths --------------------------------------------------------
# Column A of "Username" already carries the right width (21.28515625); only
# column B needs widening. "Password" needs both columns resized to match.
$wsUsername.Columns.Item(2).ColumnWidth = 40.25
$wsPassword.Columns.Item(1).ColumnWidth = 20.42
$wsPassword.Columns.Item(2).ColumnWidth = 40.25

# --- Helper: write a Text-formatted cell ---------------------------------
function Set-TextCell($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# The cell values below are written in a specific sequence so that newly
# introduced shared strings land in the same order the source workbook uses.

# Rows 1-8 (shared layout between the two sheets); column A reuses labels
# that already exist elsewhere in the workbook (FullName sheet), so order
# is not significant for those.
Set-TextCell $wsUsername "A1" "Empty"
Set-TextCell $wsUsername "B1" ""
Set-TextCell $wsUsername "A2" "OneSymbol"
Set-TextCell $wsUsername "A3" "StartWithSpace"
Set-TextCell $wsUsername "A4" "EndWithSpace"
Set-TextCell $wsUsername "A5" "MaxCharacters"
Set-TextCell $wsUsername "A6" "MinCharacters"
Set-TextCell $wsUsername "A7" "SpecialSymbols"

# New shared strings, introduced in the exact order required.
Set-TextCell $wsUsername "B3" " Username"
Set-TextCell $wsUsername "B4" "Username "
Set-TextCell $wsUsername "B5" "TheUsernameWithMoreThan32Characte"
Set-TextCell $wsUsername "B6" "Chara"
Set-TextCell $wsUsername "B7" "A!@#$%"
Set-TextCell $wsUsername "B2" "U"
Set-TextCell $wsUsername "A8" "Existing"
Set-TextCell $wsUsername "B8" "tester"
Set-TextCell $wsPassword "A9" "StartWithLatin"
Set-TextCell $wsPassword "A10" "LatinOnly"
Set-TextCell $wsUsername "A9" "NotStartWithLatin"
Set-TextCell $wsUsername "A10" "NotLatinOnly"
Set-TextCell $wsUsername "B9" "1Latin"
Set-TextCell $wsUsername "B10" "Тестер"
Set-TextCell $wsUsername "B12" "Splendor"
Set-TextCell $wsUsername "B13" "Nymgo"
Set-TextCell $wsUsername "B14" "Ahwar"
Set-TextCell $wsUsername "B15" "Callback"
Set-TextCell $wsUsername "B16" "Free"
Set-TextCell $wsUsername "B17" "Voip"
Set-TextCell $wsUsername "A12" "Reserved"
Set-TextCell $wsUsername "B11" "Tester_.tester"
Set-TextCell $wsUsername "A11" "ConsecutiveSpecial"

# Remaining "Reserved" rows on "Username" (string already registered above).
Set-TextCell $wsUsername "A13" "Reserved"
Set-TextCell $wsUsername "A14" "Reserved"
Set-TextCell $wsUsername "A15" "Reserved"
Set-TextCell $wsUsername "A16" "Reserved"
Set-TextCell $wsUsername "A17" "Reserved"

# B7 on "Username" intentionally carries no explicit cell style (unlike its
# siblings), matching the source formatting.
$wsUsername.Range("B7").NumberFormat = "General"
$wsUsername.Range("B7").Value = "A!@#$%"

# --- Mirror rows 1-8 onto "Password" (all shared strings already exist) -
Set-TextCell $wsPassword "A1" "Empty"
Set-TextCell $wsPassword "B1" ""
Set-TextCell $wsPassword "A2" "OneSymbol"
Set-TextCell $wsPassword "B2" "U"
Set-TextCell $wsPassword "A3" "StartWithSpace"
Set-TextCell $wsPassword "B3" " Username"
Set-TextCell $wsPassword "A4" "EndWithSpace"
Set-TextCell $wsPassword "B4" "Username "
Set-TextCell $wsPassword "A5" "MaxCharacters"
Set-TextCell $wsPassword "B5" "TheUsernameWithMoreThan32Characte"
Set-TextCell $wsPassword "A6" "MinCharacters"
Set-TextCell $wsPassword "B6" "Chara"
Set-TextCell $wsPassword "A7" "SpecialSymbols"
$wsPassword.Range("B7").NumberFormat = "General"
$wsPassword.Range("B7").Value = "A!@#$%"
Set-TextCell $wsPassword "A8" "Existing"
Set-TextCell $wsPassword "B8" "tester"

# --- View / active-tab bookkeeping ---------------------------------------
# The previously active sheet (FullName) loses its selection state and the
# rebuilt "Username" sheet becomes the active tab instead.
$wsFullName.Range("A1").Select() | Out-Null
$wsUsername.Activate()
$wsUsername.Range("C20").Select() | Out-Null

Write-Output "done"
